$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 6 de Octubre de 2020 a las 23:32"

# Row 4
$ws.Cells.Item(4,1).Value = "Estados Unidos"
$ws.Cells.Item(4,2).Value = 7713254
$ws.Cells.Item(4,3).Value = 33610
$ws.Cells.Item(4,4).Value = 4923013
$ws.Cells.Item(4,5).Value = 2574613
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 596
$ws.Cells.Item(4,8).Value = 215628

# Row 6
$ws.Cells.Item(6,1).Value = "Brasil"
$ws.Cells.Item(6,2).Value = 4969141
$ws.Cells.Item(6,3).Value = 28642
$ws.Cells.Item(6,4).Value = 4295302
$ws.Cells.Item(6,5).Value = 526345
$ws.Cells.Item(6,6).Value = 0
$ws.Cells.Item(6,7).Value = 721
$ws.Cells.Item(6,8).Value = 147494

# Row 13
$ws.Cells.Item(13,1).Value = "Sudafrica"
$ws.Cells.Item(13,2).Value = 683242
$ws.Cells.Item(13,3).Value = 1027
$ws.Cells.Item(13,4).Value = 616857
$ws.Cells.Item(13,5).Value = 49282
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 87
$ws.Cells.Item(13,8).Value = 17103

# Row 26
$ws.Cells.Item(26,1).Value = "Alemania"
$ws.Cells.Item(26,2).Value = 307119
$ws.Cells.Item(26,3).Value = 2462
$ws.Cells.Item(26,4).Value = 263700
$ws.Cells.Item(26,5).Value = 33784
$ws.Cells.Item(26,6).Value = 0
$ws.Cells.Item(26,7).Value = 19
$ws.Cells.Item(26,8).Value = 9635

# Row 57
$ws.Cells.Item(57,1).Value = "Barein"
$ws.Cells.Item(57,2).Value = 73476
$ws.Cells.Item(57,3).Value = 360
$ws.Cells.Item(57,4).Value = 68606
$ws.Cells.Item(57,5).Value = 4608
$ws.Cells.Item(57,6).Value = 0
$ws.Cells.Item(57,7).Value = 1
$ws.Cells.Item(57,8).Value = 262

# Row 88
$ws.Cells.Item(88,1).Value = "Costa de Marfil"
$ws.Cells.Item(88,2).Value = 19903
$ws.Cells.Item(88,3).Value = 18
$ws.Cells.Item(88,4).Value = 19539
$ws.Cells.Item(88,5).Value = 244
$ws.Cells.Item(88,6).Value = 0
$ws.Cells.Item(88,7).Value = 0
$ws.Cells.Item(88,8).Value = 120

# Row 101
$ws.Cells.Item(101,1).Value = "Namibia"
$ws.Cells.Item(101,2).Value = 11673
$ws.Cells.Item(101,3).Value = 19
$ws.Cells.Item(101,4).Value = 9611
$ws.Cells.Item(101,5).Value = 1937
$ws.Cells.Item(101,6).Value = 0
$ws.Cells.Item(101,7).Value = 0
$ws.Cells.Item(101,8).Value = 125

# Row 103
$ws.Cells.Item(103,1).Value = "Guinea"
$ws.Cells.Item(103,2).Value = 10863
$ws.Cells.Item(103,3).Value = 63
$ws.Cells.Item(103,4).Value = 10176
$ws.Cells.Item(103,5).Value = 619
$ws.Cells.Item(103,6).Value = 0
$ws.Cells.Item(103,7).Value = 1
$ws.Cells.Item(103,8).Value = 68

# Row 114
$ws.Cells.Item(114,1).Value = "Zimbabue"
$ws.Cells.Item(114,2).Value = 7915
$ws.Cells.Item(114,3).Value = 17
$ws.Cells.Item(114,4).Value = 6440
$ws.Cells.Item(114,5).Value = 1246
$ws.Cells.Item(114,6).Value = 0
$ws.Cells.Item(114,7).Value = 1
$ws.Cells.Item(114,8).Value = 229

# Row 120
$ws.Cells.Item(120,1).Value = "Malaui"
$ws.Cells.Item(120,2).Value = 5796
$ws.Cells.Item(120,3).Value = 2
$ws.Cells.Item(120,4).Value = 4549
$ws.Cells.Item(120,5).Value = 1067
$ws.Cells.Item(120,6).Value = 0
$ws.Cells.Item(120,7).Value = 0
$ws.Cells.Item(120,8).Value = 180

# Row 122
$ws.Cells.Item(122,1).Value = "Suazilandia"
$ws.Cells.Item(122,2).Value = 5598
$ws.Cells.Item(122,3).Value = 19
$ws.Cells.Item(122,4).Value = 5169
$ws.Cells.Item(122,5).Value = 317
$ws.Cells.Item(122,6).Value = 0
$ws.Cells.Item(122,7).Value = 0
$ws.Cells.Item(122,8).Value = 112

# Row 124
$ws.Cells.Item(124,1).Value = "Republica de Yibuti"
$ws.Cells.Item(124,2).Value = 5423
$ws.Cells.Item(124,3).Value = 2
$ws.Cells.Item(124,4).Value = 5353
$ws.Cells.Item(124,5).Value = 9
$ws.Cells.Item(124,6).Value = 0
$ws.Cells.Item(124,7).Value = 0
$ws.Cells.Item(124,8).Value = 61

# Row 131
$ws.Cells.Item(131,1).Value = "Ruanda"
$ws.Cells.Item(131,2).Value = 4873
$ws.Cells.Item(131,3).Value = 6
$ws.Cells.Item(131,4).Value = 3246
$ws.Cells.Item(131,5).Value = 1598
$ws.Cells.Item(131,6).Value = 0
$ws.Cells.Item(131,7).Value = 0
$ws.Cells.Item(131,8).Value = 29

# Row 134
$ws.Cells.Item(134,1).Value = "Bahamas"
$ws.Cells.Item(134,2).Value = 4559
$ws.Cells.Item(134,3).Value = 107
$ws.Cells.Item(134,4).Value = 2475
$ws.Cells.Item(134,5).Value = 1984
$ws.Cells.Item(134,6).Value = 0
$ws.Cells.Item(134,7).Value = 4
$ws.Cells.Item(134,8).Value = 100

# Row 135
$ws.Cells.Item(135,1).Value = "Siria"
$ws.Cells.Item(135,2).Value = 4457
$ws.Cells.Item(135,3).Value = 46
$ws.Cells.Item(135,4).Value = 1183
$ws.Cells.Item(135,5).Value = 3065
$ws.Cells.Item(135,6).Value = 0
$ws.Cells.Item(135,7).Value = 2
$ws.Cells.Item(135,8).Value = 209

# Row 145
$ws.Cells.Item(145,1).Value = "Mali"
$ws.Cells.Item(145,2).Value = 3195
$ws.Cells.Item(145,3).Value = 6
$ws.Cells.Item(145,4).Value = 2494
$ws.Cells.Item(145,5).Value = 570
$ws.Cells.Item(145,6).Value = 0
$ws.Cells.Item(145,7).Value = 0
$ws.Cells.Item(145,8).Value = 131

# Row 153
$ws.Cells.Item(153,1).Value = "Sierra Leona"
$ws.Cells.Item(153,2).Value = 2277
$ws.Cells.Item(153,3).Value = 8
$ws.Cells.Item(153,4).Value = 1710
$ws.Cells.Item(153,5).Value = 495
$ws.Cells.Item(153,6).Value = 0
$ws.Cells.Item(153,7).Value = 0
$ws.Cells.Item(153,8).Value = 72

# Row 156
$ws.Cells.Item(156,1).Value = "Burkina Faso"
$ws.Cells.Item(156,2).Value = 2197
$ws.Cells.Item(156,3).Value = 13
$ws.Cells.Item(156,4).Value = 1441
$ws.Cells.Item(156,5).Value = 697
$ws.Cells.Item(156,6).Value = 0
$ws.Cells.Item(156,7).Value = 0
$ws.Cells.Item(156,8).Value = 59

# Row 157
$ws.Cells.Item(157,1).Value = "Letonia"
$ws.Cells.Item(157,2).Value = 2194
$ws.Cells.Item(157,3).Value = 68
$ws.Cells.Item(157,4).Value = 1322
$ws.Cells.Item(157,5).Value = 832
$ws.Cells.Item(157,6).Value = 0
$ws.Cells.Item(157,7).Value = 1
$ws.Cells.Item(157,8).Value = 40

# Row 159
$ws.Cells.Item(159,1).Value = "Yemen"
$ws.Cells.Item(159,2).Value = 2047
$ws.Cells.Item(159,3).Value = 6
$ws.Cells.Item(159,4).Value = 1327
$ws.Cells.Item(159,5).Value = 127
$ws.Cells.Item(159,6).Value = 0
$ws.Cells.Item(159,7).Value = 1
$ws.Cells.Item(159,8).Value = 593

# Row 160
$ws.Cells.Item(160,1).Value = "Togo"
$ws.Cells.Item(160,2).Value = 1881
$ws.Cells.Item(160,3).Value = 17
$ws.Cells.Item(160,4).Value = 1410
$ws.Cells.Item(160,5).Value = 422
$ws.Cells.Item(160,6).Value = 0
$ws.Cells.Item(160,7).Value = 1
$ws.Cells.Item(160,8).Value = 49

# Row 161
$ws.Cells.Item(161,1).Value = "Republica de Chipre"
$ws.Cells.Item(161,2).Value = 1876
$ws.Cells.Item(161,3).Value = 29
$ws.Cells.Item(161,4).Value = 1369
$ws.Cells.Item(161,5).Value = 484
$ws.Cells.Item(161,6).Value = 0
$ws.Cells.Item(161,7).Value = 1
$ws.Cells.Item(161,8).Value = 23

# Row 166
$ws.Cells.Item(166,1).Value = "Republica del Chad"
$ws.Cells.Item(166,2).Value = 1238
$ws.Cells.Item(166,3).Value = 15
$ws.Cells.Item(166,4).Value = 1077
$ws.Cells.Item(166,5).Value = 73
$ws.Cells.Item(166,6).Value = 0
$ws.Cells.Item(166,7).Value = 2
$ws.Cells.Item(166,8).Value = 88

# Row 178
$ws.Cells.Item(178,1).Value = "Comoras"
$ws.Cells.Item(178,2).Value = 491
$ws.Cells.Item(178,3).Value = 4
$ws.Cells.Item(178,4).Value = 468
$ws.Cells.Item(178,5).Value = 16
$ws.Cells.Item(178,6).Value = 0
$ws.Cells.Item(178,7).Value = 0
$ws.Cells.Item(178,8).Value = 7

# Row 189
$ws.Cells.Item(189,1).Value = "Monaco"
$ws.Cells.Item(189,2).Value = 224
$ws.Cells.Item(189,3).Value = 1
$ws.Cells.Item(189,4).Value = 198
$ws.Cells.Item(189,5).Value = 24
$ws.Cells.Item(189,6).Value = 0
$ws.Cells.Item(189,7).Value = 0
$ws.Cells.Item(189,8).Value = 2

# Row 207
$ws.Cells.Item(207,1).Value = "Santa Lucia"
$ws.Cells.Item(207,2).Value = 27
$ws.Cells.Item(207,3).Value = 0
$ws.Cells.Item(207,4).Value = 27
$ws.Cells.Item(207,5).Value = 0
$ws.Cells.Item(207,6).Value = 0
$ws.Cells.Item(207,7).Value = 0
$ws.Cells.Item(207,8).Value = 0

# Row 208
$ws.Cells.Item(208,1).Value = "Nueva Caledonia"
$ws.Cells.Item(208,2).Value = 27
$ws.Cells.Item(208,3).Value = 0
$ws.Cells.Item(208,4).Value = 27
$ws.Cells.Item(208,5).Value = 0
$ws.Cells.Item(208,6).Value = 0
$ws.Cells.Item(208,7).Value = 0
$ws.Cells.Item(208,8).Value = 0
